# Applies the "Updated cryptos list" data refresh (Fri Jun 23 16:55:42 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link updates (two pairs of rows swapped position in the source feed) ---
$ws.Range("B44").Value = "Aptos"
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("B49").Value = "Elrond"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

# --- Price column updates ---
# Leading apostrophe forces text entry so values such as "1.000" or "0.9990"
# keep their exact digits/trailing zeros instead of being parsed as numbers.
$ws.Range("D2").Value = "'31.030.39"
$ws.Range("D3").Value = "'1.913.94"
$ws.Range("D4").Value = "'0.9989"
$ws.Range("D5").Value = "'246.13"
$ws.Range("D6").Value = "'0.9990"
$ws.Range("D7").Value = "'0.4965"
$ws.Range("D9").Value = "'0.06806"
$ws.Range("D10").Value = "'1.913.60"
$ws.Range("D11").Value = "'17.06"
$ws.Range("D12").Value = "'0.07309"
$ws.Range("D13").Value = "'89.83"
$ws.Range("D14").Value = "'0.6844"
$ws.Range("D15").Value = "'5.075"
$ws.Range("D16").Value = "'30.949.24"
$ws.Range("D17").Value = "'0.000008033"
$ws.Range("D18").Value = "'0.9999"
$ws.Range("D19").Value = "'13.26"
$ws.Range("D20").Value = "'2.159.71"
$ws.Range("D21").Value = "'1.005"
$ws.Range("D22").Value = "'4.885"
$ws.Range("D23").Value = "'175.78"
$ws.Range("D24").Value = "'6.076"
$ws.Range("D25").Value = "'9.344"
$ws.Range("D26").Value = "'152.30"
$ws.Range("D27").Value = "'18.12"
$ws.Range("D28").Value = "'1.951"
$ws.Range("D29").Value = "'1.439"
$ws.Range("D30").Value = "'4.341"
$ws.Range("D31").Value = "'0.08931"
$ws.Range("D32").Value = "'4.081"
$ws.Range("D33").Value = "'0.05297"
$ws.Range("D34").Value = "'0.7498"
$ws.Range("D35").Value = "'1.147"
$ws.Range("D36").Value = "'2.648"
$ws.Range("D37").Value = "'0.01932"
$ws.Range("D38").Value = "'2.730"
$ws.Range("D39").Value = "'2.213"
$ws.Range("D40").Value = "'0.9422"
$ws.Range("D41").Value = "'0.4409"
$ws.Range("D42").Value = "'5.982"
$ws.Range("D43").Value = "'105.26"
$ws.Range("D44").Value = "'7.821"
$ws.Range("D45").Value = "'1.000"
$ws.Range("D46").Value = "'0.1328"
$ws.Range("D47").Value = "'0.05848"
$ws.Range("D48").Value = "'0.3937"
$ws.Range("D49").Value = "'33.43"
$ws.Range("D50").Value = "'8.554"
$ws.Range("D51").Value = "'1.384"

# --- Volume(1h) column updates ---
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  +3.28%  "
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  +5.47%  "
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("E15").Value = "  +5.08%  "
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +4.10%  "
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("E23").Value = "  +30.68%  "
$ws.Range("E24").Value = "  +9.25%  "
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("E27").Value = "  +8.54%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  +4.25%  "
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("E33").Value = "  +6.41%  "
$ws.Range("E34").Value = "  +6.70%  "
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +18.05%  "
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("E48").Value = "  +6.77%  "
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("E50").Value = "  +4.35%  "
$ws.Range("E51").Value = "  +3.97%  "

